$wb = $excel.ActiveWorkbook

# New identifiers / values introduced by this "Generate Report for Handback" run.
$newGuid1 = "59a5d24b-0e48-4de4-987d-ad457091dcba"
$newGuid2 = "ffff66f9815e-fd99-44a1-aa5a-a27856595390"
$newHash  = "4adfb3cae53adff77fdfb5fa97466f6bbe59c6c6"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"

# ---------------------------------------------------------------------------
# Overview sheet: update the two file-name cells (A2/A3) and their hyperlinks.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2,1).Value = $newMd1
$wsOverview.Cells.Item(3,1).Value = $newMd2

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------------
# Per-locale sheets (zh-cn / de-de): update source/handoff/handback file
# names, xlf file names and the handoff/handback datetimes.
# ---------------------------------------------------------------------------
$locales = @{
    "zh-cn" = @{
        xlf2 = "$newGuid1.$newHash.zh-cn.xlf"
        xlf3 = "$newGuid1.$newHash.zh-cn.xlf"
        e2   = "2016-03-17 12:48:06"
        h2   = "2016-03-17 12:48:23"
        e3   = "2016-03-17 12:48:06"
        h3   = "2016-03-17 12:48:23"
    }
    "de-de" = @{
        xlf2 = "$newGuid1.$newHash.de-de.xlf"
        xlf3 = "$newGuid1.$newHash.de-de.xlf"
        e2   = "2016-03-17 12:48:10"
        h2   = "2016-03-17 12:48:28"
        e3   = "2016-03-17 12:48:10"
        h3   = "2016-03-17 12:48:28"
    }
}

foreach ($localeName in $locales.Keys) {
    $ws = $wb.Worksheets.Item($localeName)
    $vals = $locales[$localeName]

    # Row 2 (source file ad441f65... -> 59a5d24b...)
    $ws.Cells.Item(2,1).Value = $newMd1          # A2 Source File Name
    $ws.Cells.Item(2,4).Value = $vals.xlf2        # D2 Correspond Handoff File
    $ws.Cells.Item(2,5).Value = $vals.e2          # E2 Correspond Handoff Datetime
    $ws.Cells.Item(2,6).Value = $newMd1          # F2 Target File
    $ws.Cells.Item(2,7).Value = $vals.xlf2        # G2 Correspond Handback File
    $ws.Cells.Item(2,8).Value = $vals.h2          # H2 Correspond Handback DateTime

    # Row 3 (source file eb12e8f2... -> ffff66f9815e...)
    $ws.Cells.Item(3,1).Value = $newMd2          # A3 Source File Name
    $ws.Cells.Item(3,4).Value = $vals.xlf3        # D3 Correspond Handoff File
    $ws.Cells.Item(3,5).Value = $vals.e3          # E3 Correspond Handoff Datetime
    $ws.Cells.Item(3,6).Value = $newMd2          # F3 Target File
    $ws.Cells.Item(3,7).Value = $vals.xlf3        # G3 Correspond Handback File
    $ws.Cells.Item(3,8).Value = $vals.h3          # H3 Correspond Handback DateTime

    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        switch ($addr) {
            '$A$2' { $hl.TextToDisplay = $newMd1 }
            '$D$2' { $hl.TextToDisplay = $vals.xlf2 }
            '$F$2' { $hl.TextToDisplay = $newMd1 }
            '$G$2' { $hl.TextToDisplay = $vals.xlf2 }
            '$A$3' { $hl.TextToDisplay = $newMd2 }
            '$D$3' { $hl.TextToDisplay = $vals.xlf3 }
            '$F$3' { $hl.TextToDisplay = $newMd2 }
            '$G$3' { $hl.TextToDisplay = $vals.xlf3 }
        }
    }
}
